$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove Sheet2 and Sheet3, keep only Sheet1
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "Test myRules" table at rows 27-30
$ws.Cells.Item(27, 2).Value = "Test myRules"

$ws.Cells.Item(28, 2).Value = "aaa"
$ws.Cells.Item(28, 3).Value = "bbb"
$ws.Cells.Item(28, 4).Value = "_res_"

$ws.Cells.Item(29, 4).Value = "Result"
$ws.Cells.Item(29, 2).Value = "Arg1"
$ws.Cells.Item(29, 3).Value = "Arg2"

$ws.Cells.Item(30, 2).Value = "ccc"
$ws.Cells.Item(30, 3).Value = "yyy"

# Update selection to C30
$ws.Range("C30").Select()

Write-Host "done"
